$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44159
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = '$/bandeja 10 kilos'
$ws.Range("O2").Value = 'Provincia de Linares'
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 10

$ws.Range("D3").Value = 44159
$ws.Range("J3").Value = 260
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 9000
$ws.Range("N3").Value = '$/bandeja 10 kilos'
$ws.Range("O3").Value = 'Provincia de Linares'
$ws.Range("P3").Value = 900
$ws.Range("Q3").Value = 10

$ws.Range("D4").Value = 44159
$ws.Range("J4").Value = 320
$ws.Range("K4").Value = 7000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 7000
$ws.Range("N4").Value = '$/bandeja 10 kilos'
$ws.Range("O4").Value = 'Provincia de Linares'
$ws.Range("P4").Value = 700
$ws.Range("Q4").Value = 10

$ws.Range("D5").Value = 44169
$ws.Range("H5").Value = 'Verde'
$ws.Range("I5").Value = 'Banquete'
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 10000
$ws.Range("N5").Value = '$/bandeja 10 kilos'
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = 10

$ws.Range("D6").Value = 44169
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 95
$ws.Range("K6").Value = 7500
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = 7500
$ws.Range("N6").Value = '$/bandeja 10 kilos'
$ws.Range("P6").Value = 750

$ws.Range("D7").Value = 44169
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 110
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = 6500
$ws.Range("N7").Value = '$/bandeja 10 kilos'
$ws.Range("P7").Value = 650

$ws.Range("D8").Value = 44161
$ws.Range("I8").Value = 'Banquete'
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9500
$ws.Range("M8").Value = 9269
$ws.Range("O8").Value = 'Provincia de Linares'
$ws.Range("P8").Value = 927

$ws.Range("D9").Value = 44161
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 7500
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7778
$ws.Range("N9").Value = '$/caja 10 kilos'
$ws.Range("O9").Value = 'Provincia de Linares'
$ws.Range("P9").Value = 778

$ws.Range("D10").Value = 44161
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 6000
$ws.Range("L10").Value = 6500
$ws.Range("M10").Value = 6300
$ws.Range("N10").Value = '$/caja 10 kilos'
$ws.Range("O10").Value = 'Provincia de Linares'
$ws.Range("P10").Value = 630

$ws.Range("D11").Value = 44453
$ws.Range("H11").Value = 'Sin especificar'
$ws.Range("J11").Value = 55
$ws.Range("K11").Value = 2600
$ws.Range("L11").Value = 2600
$ws.Range("M11").Value = 2600
$ws.Range("N11").Value = '$/kilo'
$ws.Range("O11").Value = 'Provincia de Linares'
$ws.Range("P11").Value = 2600
$ws.Range("Q11").Value = 1

$ws.Range("D12").Value = 44168
$ws.Range("J12").Value = 105
$ws.Range("K12").Value = 9000
$ws.Range("M12").Value = 9429
$ws.Range("N12").Value = '$/caja 10 kilos'
$ws.Range("P12").Value = 943

$ws.Range("D13").Value = 44168
$ws.Range("J13").Value = 290
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 7500
$ws.Range("M13").Value = 7241
$ws.Range("N13").Value = '$/caja 10 kilos'
$ws.Range("P13").Value = 724

$ws.Range("D14").Value = 44168
$ws.Range("J14").Value = 360
$ws.Range("K14").Value = 6000
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6278
$ws.Range("N14").Value = '$/caja 10 kilos'
$ws.Range("P14").Value = 628

$ws.Range("D15").Value = 44160
$ws.Range("J15").Value = 400
$ws.Range("L15").Value = 9500
$ws.Range("M15").Value = 9275
$ws.Range("N15").Value = '$/bandeja 8 kilos'
$ws.Range("O15").Value = 'Región Metropolitana'
$ws.Range("P15").Value = 1159
$ws.Range("Q15").Value = 8

$ws.Range("D16").Value = 44160
$ws.Range("J16").Value = 440
$ws.Range("K16").Value = 7500
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 7784
$ws.Range("N16").Value = '$/bandeja 8 kilos'
$ws.Range("O16").Value = 'Región Metropolitana'
$ws.Range("P16").Value = 973
$ws.Range("Q16").Value = 8

$ws.Range("D17").Value = 44160
$ws.Range("J17").Value = 305
$ws.Range("M17").Value = 6270
$ws.Range("N17").Value = '$/bandeja 8 kilos'
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 784
$ws.Range("Q17").Value = 8

$ws.Range("D18").Value = 44175
$ws.Range("J18").Value = 120
$ws.Range("N18").Value = '$/caja 10 kilos'
$ws.Range("O18").Value = 'Provincia de Linares'

$ws.Range("D19").Value = 44175
$ws.Range("J19").Value = 140
$ws.Range("N19").Value = '$/caja 10 kilos'
$ws.Range("O19").Value = 'Provincia de Linares'

$ws.Range("D20").Value = 44162
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 10000
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("P20").Value = 1000

$ws.Range("D21").Value = 44162
$ws.Range("J21").Value = 220
$ws.Range("K21").Value = 8500
$ws.Range("L21").Value = 8500
$ws.Range("M21").Value = 8500
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 850

$ws.Range("D22").Value = 44162
$ws.Range("J22").Value = 260
$ws.Range("K22").Value = 7500
$ws.Range("L22").Value = 7500
$ws.Range("M22").Value = 7500
$ws.Range("O22").Value = 'Región Metropolitana'
$ws.Range("P22").Value = 750

$ws.Range("D23").Value = 44176
$ws.Range("J23").Value = 80
$ws.Range("N23").Value = '$/bandeja 10 kilos'
$ws.Range("O23").Value = 'Región Metropolitana'

$ws.Range("D24").Value = 44176
$ws.Range("J24").Value = 100
$ws.Range("N24").Value = '$/bandeja 10 kilos'
$ws.Range("O24").Value = 'Región Metropolitana'
